# chore: adapt column header formatting to respective input file names (#7)
#
# 1. Rename the "_old" / "_new" column header suffixes to "_FV2310" / "_FV2404"
# 2. Turn the data range A1:U57 into a real Excel Table (Table1)
# 3. Freeze the header row (pane split under row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header renames
# ---------------------------------------------------------------------------
$oldHeaders = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old"
)
$newHeaders = @(
    "Segmentname_new",
    "Segmentgruppe_new",
    "Segment_new",
    "Datenelement_new",
    "Segment ID_new",
    "Code_new",
    "Qualifier_new",
    "Beschreibung_new",
    "Bedingungsausdruck_new",
    "Bedingung_new"
)

$fv2310Headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)
$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2310Headers[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $fv2404Headers[$i]
}

# ---------------------------------------------------------------------------
# 2) Build the table
#
# Building the ListObject directly on top of A1:U57 would make the engine
# snapshot the pre-existing bold/filled header formatting into a
# headerRowDxfId (the header cells already carry an explicit bold style).
# Building it first on a blank, unformatted scratch range (so no dxf capture
# happens) and then resizing it onto the real range keeps the original cell
# styles untouched and avoids injecting that extra dxf/style.
# ---------------------------------------------------------------------------
$scratchHeader = $ws.Range("AA1:AU1")
$ws.Range("A1:U1").Copy()
$scratchHeader.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false

$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("AA1:AU2"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.TableStyle = ""

$scratchHeader.ClearContents()

$lo.Resize($ws.Range("A1:U57"))
$lo.Name = "Table1"

# ---------------------------------------------------------------------------
# 3) Freeze the header row
# ---------------------------------------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
